# edit.ps1 - apply the authored changes to the presentation:
#   1. Slide 5's table: switch its table style to the built-in
#      {08526BC9-1CF3-45F1-B89D-484876083323} style.
#   2. The deck's theme colour scheme: swap the "Integral / Red Violet"
#      palette for the default "Office" palette (dk1/lt1/dk2/lt2/accent1-6/
#      hlink/folHlink), as it ends up after the authored edit.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 --------------------------------------------
$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{08526BC9-1CF3-45F1-B89D-484876083323}")
    }
}

# --- 2. Theme colour scheme -------------------------------------------------
# New (target) "Office" theme colours, in the standard ppColorSchemeIndex
# order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeColors = @(
    0x000000,
    0xFFFFFF,
    0x44546A,
    0xE7E6E6,
    0x5B9BD5,
    0xED7D31,
    0xA5A5A5,
    0xFFC000,
    0x4472C4,
    0x70AD47,
    0x0563C1,
    0x954F72
)

function ToCOMRGB($hexColor) {
    # PowerPoint RGB values are stored as 0x00BBGGRR, not 0x00RRGGBB.
    $r = ($hexColor -band 0xFF0000) -shr 16
    $g = ($hexColor -band 0x00FF00) -shr 8
    $b = ($hexColor -band 0x0000FF)
    return $b * 65536 + $g * 256 + $r
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = ToCOMRGB($officeColors[$i - 1])
}
